# Regenerate the localization-status report: both zh-cn and de-de are now
# handed back and in sync with en-US, and the latest handback round for
# both locales was re-generated on 2017-02-28 (with the stale "handback not
# latest" warning cleared now that everything is current).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: summary status column for each locale ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $statusText
$ovw.Range("F2").Value = $statusText
$ovw.Columns.Item(5).ColumnWidth = 29.1666666667
$ovw.Columns.Item(6).ColumnWidth = 29.1666666667

# --- zh-cn detail sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $statusText
$zh.Range("L2").Value = "2017-02-28 07:08:33"
$zh.Range("M2").Value = "TestHandback_201702280308"
$zh.Range("R2").Value = ""
$zh.Columns.Item(3).ColumnWidth = 29.1666666667
$zh.Columns.Item(18).ColumnWidth = 12.8333333333

# --- de-de detail sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $statusText
$de.Range("L2").Value = "2017-02-28 07:08:55"
$de.Range("M2").Value = "TestHandback_201702280308"
$de.Range("R2").Value = ""
$de.Columns.Item(3).ColumnWidth = 29.1666666667
$de.Columns.Item(18).ColumnWidth = 12.8333333333
